{"js": "// Replace the two-digit-divided-by-one-digit answer strings in the\n// worksheet table with a new set of problems/answers, cell by cell,\n// preserving all existing formatting (fonts, sizes, paragraph props).\n//\n// The old values are applied/verified positionally (row-major, left to\n// right) rather than via global find/replace, because some old values\n// repeat (e.g. \"22\u00f78=2, 6\" appears twice) but map to different new\n// values depending on position.\nconst oldNewPairs = [\n  [\"82\u00f77=11, 5\", \"62\u00f79=6, 8\"],\n  [\"29\u00f72=14, 1\", \"59\u00f73=19, 2\"],\n  [\"65\u00f77=9, 2\", \"64\u00f72=32, 0\"],\n  [\"48\u00f73=16, 0\", \"89\u00f77=12, 5\"],\n  [\"38\u00f79=4, 2\", \"51\u00f77=7, 2\"],\n  [\"37\u00f74=9, 1\", \"30\u00f78=3, 6\"],\n  [\"46\u00f76=7, 4\", \"48\u00f74=12, 0\"],\n  [\"94\u00f72=47, 0\", \"33\u00f74=8, 1\"],\n  [\"56\u00f77=8, 0\", \"20\u00f76=3, 2\"],\n  [\"64\u00f73=21, 1\", \"90\u00f76=15, 0\"],\n  [\"10\u00f74=2, 2\", \"62\u00f75=12, 2\"],\n  [\"53\u00f72=26, 1\", \"42\u00f79=4, 6\"],\n  [\"85\u00f72=42, 1\", \"75\u00f73=25, 0\"],\n  [\"99\u00f78=12, 3\", \"66\u00f78=8, 2\"],\n  [\"22\u00f78=2, 6\", \"22\u00f77=3, 1\"],\n  [\"39\u00f73=13, 0\", \"22\u00f75=4, 2\"],\n  [\"34\u00f75=6, 4\", \"81\u00f79=9, 0\"],\n  [\"22\u00f78=2, 6\", \"36\u00f79=4, 0\"],\n  [\"60\u00f76=10, 0\", \"78\u00f74=19, 2\"],\n  [\"75\u00f74=18, 3\", \"74\u00f79=8, 2\"],\n  [\"63\u00f79=7, 0\", \"91\u00f76=15, 1\"],\n  [\"69\u00f78=8, 5\", \"53\u00f75=10, 3\"],\n  [\"18\u00f79=2, 0\", \"96\u00f73=32, 0\"],\n  [\"53\u00f77=7, 4\", \"49\u00f74=12, 1\"],\n  [\"63\u00f76=10, 3\", \"99\u00f72=49, 1\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst values = table.values;\n\n// Walk the table in row-major order, consuming one (old, new) pair for\n// every non-empty cell encountered (the sheet also has blank spacer\n// rows between the rows of problems, which we skip over).\nlet pairIndex = 0;\nfor (let r = 0; r < values.length; r++) {\n  for (let c = 0; c < values[r].length; c++) {\n    const cellText = values[r][c];\n    if (cellText === \"\") {\n      continue;\n    }\n    const [expectedOld, newText] = oldNewPairs[pairIndex];\n    if (cellText !== expectedOld) {\n      throw new Error(\n        `Unexpected cell text at row ${r}, col ${c}: ` +\n        `expected \"${expectedOld}\" but found \"${cellText}\"`\n      );\n    }\n    values[r][c] = newText;\n    pairIndex++;\n  }\n}\n\ntable.values = values;\nawait context.sync();\n", "ps1": "# Replace the two-digit-divided-by-one-digit answer strings in the\n# worksheet table with a new set of problems/answers, cell by cell,\n# preserving all existing formatting (fonts, sizes, paragraph props).\n#\n# The old values are applied/verified positionally (row-major, left to\n# right) rather than via global find/replace, because some old values\n# repeat (e.g. \"22\u00f78=2, 6\" appears twice) but map to different new\n# values depending on position.\n$d = $word.ActiveDocument\n\n$oldNewPairs = @(\n    , @(\"82\u00f77=11, 5\", \"62\u00f79=6, 8\")\n    , @(\"29\u00f72=14, 1\", \"59\u00f73=19, 2\")\n    , @(\"65\u00f77=9, 2\", \"64\u00f72=32, 0\")\n    , @(\"48\u00f73=16, 0\", \"89\u00f77=12, 5\")\n    , @(\"38\u00f79=4, 2\", \"51\u00f77=7, 2\")\n    , @(\"37\u00f74=9, 1\", \"30\u00f78=3, 6\")\n    , @(\"46\u00f76=7, 4\", \"48\u00f74=12, 0\")\n    , @(\"94\u00f72=47, 0\", \"33\u00f74=8, 1\")\n    , @(\"56\u00f77=8, 0\", \"20\u00f76=3, 2\")\n    , @(\"64\u00f73=21, 1\", \"90\u00f76=15, 0\")\n    , @(\"10\u00f74=2, 2\", \"62\u00f75=12, 2\")\n    , @(\"53\u00f72=26, 1\", \"42\u00f79=4, 6\")\n    , @(\"85\u00f72=42, 1\", \"75\u00f73=25, 0\")\n    , @(\"99\u00f78=12, 3\", \"66\u00f78=8, 2\")\n    , @(\"22\u00f78=2, 6\", \"22\u00f77=3, 1\")\n    , @(\"39\u00f73=13, 0\", \"22\u00f75=4, 2\")\n    , @(\"34\u00f75=6, 4\", \"81\u00f79=9, 0\")\n    , @(\"22\u00f78=2, 6\", \"36\u00f79=4, 0\")\n    , @(\"60\u00f76=10, 0\", \"78\u00f74=19, 2\")\n    , @(\"75\u00f74=18, 3\", \"74\u00f79=8, 2\")\n    , @(\"63\u00f79=7, 0\", \"91\u00f76=15, 1\")\n    , @(\"69\u00f78=8, 5\", \"53\u00f75=10, 3\")\n    , @(\"18\u00f79=2, 0\", \"96\u00f73=32, 0\")\n    , @(\"53\u00f77=7, 4\", \"49\u00f74=12, 1\")\n    , @(\"63\u00f76=10, 3\", \"99\u00f72=49, 1\")\n)\n\n$table = $d.Tables.Item(1)\n$pairIndex = 0\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    for ($c = 1; $c -le $table.Columns.Count; $c++) {\n        $cell = $table.Cell($r, $c)\n        # Cell.Range.Text includes the trailing cell-mark (CR + BEL);\n        # strip it to get just the visible text.\n        $cellText = $cell.Range.Text.TrimEnd([char]7).TrimEnd([char]13)\n        if ($cellText -eq \"\") {\n            continue\n        }\n        $expectedOld = $oldNewPairs[$pairIndex][0]\n        $newText = $oldNewPairs[$pairIndex][1]\n        if ($cellText -ne $expectedOld) {\n            throw \"Unexpected cell text at row $r, col ${c}: expected '$expectedOld' but found '$cellText'\"\n        }\n        $cell.Range.Text = $newText\n        $pairIndex++\n    }\n}\n"}
